$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared string values (sharedStrings.xml) ---
$ws.Range("A1").Value = "indexador_dotacao"
$ws.Range("A2").Value = "2018.16.2100.3.3.90.30.00"
$ws.Range("B2").Value = "categoria "
$ws.Range("C2").Value = "categoria desc "
$ws.Range("A3").Value = "2018.16.2100.3.3.90.30.99"
$ws.Range("B3").Value = "outra categoria "
$ws.Range("C3").Value = "outra categoria desc "

# --- Add two new defined names (workbook.xml), sheet-scoped like the existing ones ---
$rng = $ws.Range("A1:C341")
$ws.Names.Add("_xlnm._FilterDatabase_0_0", $rng)
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0", $rng)

# --- Move the selection from D1 to A4 (sheet1.xml sheetView/selection) ---
$null = $ws.Range("A4").Select()
